# Remove the trailing "Ver no Jupiter..." / "(c) 2020 ..." footer block,
# along with the blank paragraph that precedes them, from the end of the
# Bibliografia section. The paragraph containing
# "Rio de Janeiro: Elsevier Editora, 2007." and the blank paragraph that
# follows the removed block must remain untouched.

$d = $word.ActiveDocument

$jupiterText = "Ver no Jupiter Salvar em pdf Salvar em docx"
$copyrightText = "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"

# Walk paragraphs from the end towards the start so deleting a paragraph
# doesn't shift the indices of paragraphs we still need to inspect.
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text.Trim()

    if ($text -eq $copyrightText -or $text -eq $jupiterText) {
        # Delete this paragraph entirely (text + paragraph mark).
        $para.Range.Delete()
    }
}

# Now remove the now-orphaned blank paragraph that used to sit between
# "Rio de Janeiro: Elsevier Editora, 2007." and "Ver no Jupiter ...".
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text.Trim()

    if ($text -eq "") {
        $prevText = ""
        if ($i -gt 1) {
            $prevText = $d.Paragraphs.Item($i - 1).Range.Text.Trim()
        }
        if ($prevText -eq "Rio de Janeiro: Elsevier Editora, 2007.") {
            $para.Range.Delete()
            break
        }
    }
}
